$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$scratch = $ws.Range("ZZ9999")

$ws.Range('D2').Value = '27.706.01'
$ws.Range('E2').Value = '  -0.33%  '
$ws.Range('D3').Value = '1.846.83'
$ws.Range('E3').Value = '  -1.05%  '
$scratch.Value = "'1.009"
$scratch.Copy()
$ws.Range('D4').PasteSpecial(-4163)
$ws.Range('E4').Value = '  -3.02%  '
$scratch.Value = "'319.33"
$scratch.Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E5').Value = '  -1.68%  '
$ws.Range('E6').Value = '  -2.60%  '
$scratch.Value = "'0.4312"
$scratch.Copy()
$ws.Range('D7').PasteSpecial(-4163)
$ws.Range('E7').Value = '  -2.63%  '
$scratch.Value = "'0.3748"
$scratch.Copy()
$ws.Range('D8').PasteSpecial(-4163)
$ws.Range('E8').Value = '  -1.52%  '
$scratch.Value = "'0.07357"
$scratch.Copy()
$ws.Range('D9').PasteSpecial(-4163)
$ws.Range('E9').Value = '  -1.67%  '
$scratch.Value = "'0.8812"
$scratch.Copy()
$ws.Range('D10').PasteSpecial(-4163)
$ws.Range('E10').Value = '  -0.58%  '
$scratch.Value = "'21.61"
$scratch.Copy()
$ws.Range('D11').PasteSpecial(-4163)
$ws.Range('E11').Value = '  -0.60%  '
$ws.Range('D12').Value = '1.848.75'
$ws.Range('E12').Value = '  -1.29%  '
$scratch.Value = "'6.738"
$scratch.Copy()
$ws.Range('D13').PasteSpecial(-4163)
$ws.Range('E13').Value = '  -0.37%  '
$scratch.Value = "'5.457"
$scratch.Copy()
$ws.Range('D14').PasteSpecial(-4163)
$ws.Range('E14').Value = '  -1.98%  '
$scratch.Value = "'0.07141"
$scratch.Copy()
$ws.Range('D15').PasteSpecial(-4163)
$ws.Range('E15').Value = '  -1.36%  '
$scratch.Value = "'87.95"
$scratch.Copy()
$ws.Range('D16').PasteSpecial(-4163)
$ws.Range('E16').Value = '  +4.75%  '
$scratch.Value = "'1.013"
$scratch.Copy()
$ws.Range('D17').PasteSpecial(-4163)
$ws.Range('E17').Value = '  -2.84%  '
$scratch.Value = "'0.000009001"
$scratch.Copy()
$ws.Range('D18').PasteSpecial(-4163)
$ws.Range('E18').Value = '  -1.73%  '
$scratch.Value = "'1.010"
$scratch.Copy()
$ws.Range('D19').PasteSpecial(-4163)
$ws.Range('E19').Value = '  -2.62%  '
$ws.Range('E20').Value = '  -0.58%  '
$ws.Range('D21').Value = '27.706.34'
$ws.Range('E21').Value = '  -0.37%  '
$scratch.Value = "'5.253"
$scratch.Copy()
$ws.Range('D22').PasteSpecial(-4163)
$ws.Range('E22').Value = '  -1.42%  '
$ws.Range('E23').Value = '  -2.00%  '
$ws.Range('D24').Value = '2.075.62'
$ws.Range('E24').Value = '  -1.43%  '
$scratch.Value = "'2.020"
$scratch.Copy()
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('E25').Value = '  +1.61%  '
$scratch.Value = "'155.60"
$scratch.Copy()
$ws.Range('D26').PasteSpecial(-4163)
$ws.Range('E26').Value = '  -2.03%  '
$scratch.Value = "'18.60"
$scratch.Copy()
$ws.Range('D27').PasteSpecial(-4163)
$ws.Range('E27').Value = '  -1.65%  '
$scratch.Value = "'2.136"
$scratch.Copy()
$ws.Range('D28').PasteSpecial(-4163)
$ws.Range('E28').Value = '  +7.44%  '
$scratch.Value = "'5.400"
$scratch.Copy()
$ws.Range('D29').PasteSpecial(-4163)
$ws.Range('E29').Value = '  +1.21%  '
$scratch.Value = "'120.09"
$scratch.Copy()
$ws.Range('D30').PasteSpecial(-4163)
$ws.Range('E30').Value = '  +1.84%  '
$scratch.Value = "'0.08926"
$scratch.Copy()
$ws.Range('D31').PasteSpecial(-4163)
$ws.Range('E31').Value = '  -1.75%  '
$scratch.Value = "'1.234"
$scratch.Copy()
$ws.Range('D32').PasteSpecial(-4163)
$ws.Range('E32').Value = '  +1.40%  '
$ws.Range('E33').Value = '  +0.07%  '
$scratch.Value = "'4.564"
$scratch.Copy()
$ws.Range('D34').PasteSpecial(-4163)
$ws.Range('E34').Value = '  -0.40%  '
$scratch.Value = "'2.912"
$scratch.Copy()
$ws.Range('D35').PasteSpecial(-4163)
$ws.Range('E35').Value = '  -6.77%  '
$ws.Range('E36').Value = '  -2.76%  '
$scratch.Value = "'1.137"
$scratch.Copy()
$ws.Range('D37').PasteSpecial(-4163)
$ws.Range('E37').Value = '  -1.49%  '
$scratch.Value = "'0.05335"
$scratch.Copy()
$ws.Range('D38').PasteSpecial(-4163)
$ws.Range('E38').Value = '  -0.38%  '
$scratch.Value = "'0.01974"
$scratch.Copy()
$ws.Range('D39').PasteSpecial(-4163)
$ws.Range('E39').Value = '  -1.12%  '
$scratch.Value = "'7.247"
$scratch.Copy()
$ws.Range('D40').PasteSpecial(-4163)
$ws.Range('E40').Value = '  +4.45%  '
$scratch.Value = "'2.869"
$scratch.Copy()
$ws.Range('D41').PasteSpecial(-4163)
$ws.Range('E41').Value = '  +0.12%  '
$scratch.Value = "'0.5167"
$scratch.Copy()
$ws.Range('D42').PasteSpecial(-4163)
$ws.Range('E42').Value = '  -0.85%  '
$ws.Range('E43').Value = '  -1.06%  '
$scratch.Value = "'8.941"
$scratch.Copy()
$ws.Range('D44').PasteSpecial(-4163)
$ws.Range('E44').Value = '  +2.93%  '
$scratch.Value = "'110.12"
$scratch.Copy()
$ws.Range('D45').PasteSpecial(-4163)
$ws.Range('E45').Value = '  +0.31%  '
$scratch.Value = "'10.70"
$scratch.Copy()
$ws.Range('D46').PasteSpecial(-4163)
$ws.Range('E46').Value = '  +0.18%  '
$scratch.Value = "'0.4736"
$scratch.Copy()
$ws.Range('D47').PasteSpecial(-4163)
$ws.Range('E47').Value = '  +0.32%  '
$scratch.Value = "'0.06503"
$scratch.Copy()
$ws.Range('D48').PasteSpecial(-4163)
$ws.Range('E48').Value = '  +0.36%  '
$scratch.Value = "'1.700"
$scratch.Copy()
$ws.Range('D49').PasteSpecial(-4163)
$ws.Range('E49').Value = '  -1.59%  '
$ws.Range('E50').Value = '  -2.77%  '
$scratch.Value = "'1.895"
$scratch.Copy()
$ws.Range('D51').PasteSpecial(-4163)
$ws.Range('E51').Value = '  -0.93%  '
$scratch.Clear()
